# Add three new backlog items to the bottom of the Table1 list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: finish off the previously-empty "What" cell for item #4
$ws.Range("B5").Value = "Add Metadata files path as configurable path"

# New row 6: item #5
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Add SKIP UI/Page generation to build and only build UI code"

# New row 7: item #6
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Update navigation component of UI"

# Grow the Table1 ListObject (and its AutoFilter) to include the new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D7"))

# Widen column B so the longer text fits
$ws.Columns.Item(2).ColumnWidth = 49.5

# Mirror the author's final selection state
$ws.Range("D6:D7").Select()
